$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-39 (A: index, B, C, D)
$data = @(
    @(0, 3.707140224301534, 3.178037295964893, 0.3141922284839216),
    @(1, 5.769324542341939, 3.514645759822572, 0.228229366786155),
    @(2, 7.62411047044329, 9.98737337772222, 0.3564254466116638),
    @(3, 9.43712470116929, 4.526633084058997, 0.1698233758727088),
    @(4, 9.789298309195297, 10.45814216566142, 0.3797137767678765),
    @(5, 13.34810213934554, 2.082284836937426, 0.2786643713207221),
    @(6, 13.69850290193687, 8.222285311106274, 0.4568447624185565),
    @(7, 15.35732177268107, 8.222141808501618, 0.2768102617653607),
    @(8, 20.00453855248269, 9.480297720733683, 0.334879040270145),
    @(9, 26.46406827117856, 6.263165106537548, 0.4644998202591301),
    @(10, 28.68562946610409, 6.506682668832935, 0.2192083281815637),
    @(11, 29.78852021004996, 8.691168402951606, 0.4319581193448346),
    @(12, 30.6638744015881, 2.854514251357171, 0.08354427686892525),
    @(13, 32.80963041655549, 4.162630145626276, 0.5161727024809879),
    @(14, 36.95127302187949, 13.96310513455267, 0.3695092066035308),
    @(15, 39.8950950767116, 3.968872534344617, 0.2875625288617086),
    @(16, 45.78756001447373, 7.151446291230044, 0.2463038652876009),
    @(17, 46.15947409014645, 9.582716971428685, 0.254337097059642),
    @(18, 47.21662386988176, 3.892015393812848, 0.1596919743314569),
    @(19, 47.47509989295278, 5.9257383234276, 0.2068762850183046),
    @(20, 48.75873895281444, 5.391650038909797, 0.445702046923163),
    @(21, 53.51416894627793, 3.053275970643026, 0.1874779464143486),
    @(22, 57.44856667250019, 4.357603152176637, 0.06687991237364961),
    @(23, 59.10289925521387, 8.191369814116905, 0.3513996493502551),
    @(24, 64.45474118187454, 8.088874663895576, 0.4214097781080564),
    @(25, 65.60072671935656, 12.31845570030995, 0.4874573600337201),
    @(26, 75.20973427110188, 8.141894827849525, 0.3164131369174463),
    @(27, 84.99821336541275, 11.82686273112039, 0.6022214615835444),
    @(28, 85.48502343959228, 5.876111935471082, 0.3209055237613112),
    @(29, 85.64338473698835, 10.11446858351495, 0.4422328438256813),
    @(30, 87.69507773187695, 6.634584531426813, 0.1709929670522698),
    @(31, 87.90016609094107, 7.97895600623199, 0.2275441664472891),
    @(32, 88.15776302861563, 9.138825674050988, 0.3599001097425593),
    @(33, 90.84167439540904, 7.705029194634697, 0.399091601216635),
    @(34, 94.04380190590184, 4.278583697672139, 0.3781578102582177),
    @(35, 94.44914731532552, 8.188989842533836, 0.371252639607788),
    @(36, 97.67280456037174, 5.334572805196472, 0.4935353404420858),
    @(37, 98.46310987495943, 7.782947630754367, 0.4010391167315959)
)

# Update existing rows 2..39 with new B, C, D values (A column stays the same index values)
$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

# Delete the now-extra rows 40-44 (the sheet previously had data through row 44)
# xlShiftUp = -4162
$ws.Range("A40:D44").Delete(-4162)
